# Sprint 40 test case report - add Day 8 (row 44-47 block) summary numbers
# and update the saved sheet view/selection to reflect the new active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Day 8" block totals that were previously blank.
$ws.Range("C45").Value = 2879   # Total testcase Written
$ws.Range("C46").Value = 1426   # Total Execution
$ws.Range("C47").Value = 823    # Total Review

# Move the selection/view down to the newly filled-in row, matching the
# saved workbook view (topLeftCell A33, active cell C47).
[void]$ws.Range("C47").Select()
